$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force plain-numeric-looking price text cells to Text format so Excel keeps
# them as literal strings (matching trailing zeros / precision) instead of
# silently converting them to floating point numbers.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D2").Value = "63.252.94"
$ws.Range("E2").Value = "  +2.90%  "
$ws.Range("D3").Value = "3.488.26"
$ws.Range("E3").Value = "  +2.86%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "584.85"
$ws.Range("E5").Value = "  +1.63%  "
$ws.Range("D6").Value = "148.13"
$ws.Range("E6").Value = "  +5.39%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "0.480"
$ws.Range("E8").Value = "  +1.27%  "
$ws.Range("E9").Value = "  +0.36%  "
$ws.Range("E10").Value = "  +3.09%  "
$ws.Range("D11").Value = "0.398"
$ws.Range("E11").Value = "  +2.59%  "
$ws.Range("D12").Value = "4.085.75"
$ws.Range("E12").Value = "  +2.88%  "
$ws.Range("D13").Value = "29.76"
$ws.Range("E13").Value = "  +5.05%  "
$ws.Range("E14").Value = "  -0.13%  "
$ws.Range("D15").Value = "3.486.89"
$ws.Range("E15").Value = "  +2.75%  "
$ws.Range("D16").Value = "0.0000173"
$ws.Range("E16").Value = "  +2.06%  "
$ws.Range("D17").Value = "63.276.38"
$ws.Range("E17").Value = "  +2.95%  "
$ws.Range("D18").Value = "6.30"
$ws.Range("E18").Value = "  +2.38%  "
$ws.Range("D19").Value = "14.37"
$ws.Range("E19").Value = "  +5.16%  "
$ws.Range("D20").Value = "9.38"
$ws.Range("E20").Value = "  +4.51%  "
$ws.Range("D21").Value = "389.78"
$ws.Range("E21").Value = "  +0.11%  "
$ws.Range("E22").Value = "  +1.95%  "
$ws.Range("D23").Value = "75.22"
$ws.Range("E23").Value = "  -0.01%  "
$ws.Range("E24").Value = "  -0.06%  "
$ws.Range("D25").Value = "0.0000118"
$ws.Range("E25").Value = "  +5.23%  "
$ws.Range("D26").Value = "3.630.34"
$ws.Range("E26").Value = "  +2.90%  "
$ws.Range("E27").Value = "  -4.56%  "
$ws.Range("D28").Value = "7.80"
$ws.Range("E28").Value = "  +6.77%  "
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("D30").Value = "8.30"
$ws.Range("E30").Value = "  +3.54%  "
$ws.Range("D31").Value = "1.46"
$ws.Range("E31").Value = "  +5.91%  "
$ws.Range("E32").Value = "  +0.28%  "
$ws.Range("E33").Value = "  -0.03%  "
$ws.Range("D34").Value = "23.86"
$ws.Range("E34").Value = "  +1.97%  "
$ws.Range("D35").Value = "5.36"
$ws.Range("E35").Value = "  +6.58%  "
$ws.Range("B36").Value = "EnergySwap"
$ws.Range("C36").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D36").Value = "32.39"
$ws.Range("E36").Value = "  +25.58%  "
$ws.Range("B37").Value = "Aptos"
$ws.Range("C37").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D37").Value = "7.14"
$ws.Range("E37").Value = "  +3.16%  "
$ws.Range("D38").Value = "171.55"
$ws.Range("E38").Value = "  +2.42%  "
$ws.Range("D39").Value = "1.58"
$ws.Range("E39").Value = "  +6.94%  "
$ws.Range("D40").Value = "3.524.27"
$ws.Range("E40").Value = "  +2.85%  "
$ws.Range("D41").Value = "0.0770"
$ws.Range("E41").Value = "  +0.25%  "
$ws.Range("D42").Value = "0.808"
$ws.Range("E42").Value = "  +3.97%  "
$ws.Range("B43").Value = "Filecoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D43").Value = "4.50"
$ws.Range("E43").Value = "  +1.32%  "
$ws.Range("B44").Value = "Stacks"
$ws.Range("C44").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D44").Value = "1.73"
$ws.Range("E44").Value = "  +4.25%  "
$ws.Range("D45").Value = "42.46"
$ws.Range("E45").Value = "  +0.19%  "
$ws.Range("E46").Value = "  +7.49%  "
$ws.Range("D47").Value = "2.629.56"
$ws.Range("E47").Value = "  +7.34%  "
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").Value = "23.63"
$ws.Range("E48").Value = "  +3.70%  "
$ws.Range("B49").Value = "dogwifhat"
$ws.Range("C49").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D49").Value = "2.30"
$ws.Range("E49").Value = "  +12.78%  "
$ws.Range("D50").Value = "6.78"
$ws.Range("E50").Value = "  +1.24%  "
$ws.Range("E51").Value = "  +3.50%  "
